$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("D13").Value = 44467
$ws.Range("M13").Value = 200

# Row 14
$ws.Range("D14").Value = 44434
$ws.Range("M14").Value = 100

# Row 15
$ws.Range("D15").Value = 44441

# Row 16
$ws.Range("D16").Value = 44407
$ws.Range("M16").Value = 160
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 20500
$ws.Range("S16").Value = 1025

# Row 17
$ws.Range("D17").Value = 44336
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 19500
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 19750
$ws.Range("S17").Value = 988

# Row 18
$ws.Range("D18").Value = 44442
$ws.Range("M18").Value = 140

# Row 19
$ws.Range("D19").Value = 44448
$ws.Range("M19").Value = 100

# Row 20
$ws.Range("D20").Value = 44420
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 20500
$ws.Range("S20").Value = 1025

# Row 21
$ws.Range("D21").Value = 44350
$ws.Range("M21").Value = 160
$ws.Range("N21").Value = 19000
$ws.Range("O21").Value = 20000
$ws.Range("P21").Value = 19500
$ws.Range("S21").Value = 975

# Row 22
$ws.Range("D22").Value = 44418

# Row 23
$ws.Range("D23").Value = 44427
$ws.Range("M23").Value = 200

# Row 24
$ws.Range("D24").Value = 44466
$ws.Range("N24").Value = 20000
$ws.Range("O24").Value = 21000
$ws.Range("P24").Value = 20500
$ws.Range("S24").Value = 1025

# Row 25
$ws.Range("D25").Value = 44343
$ws.Range("N25").Value = 19500
$ws.Range("O25").Value = 20000
$ws.Range("P25").Value = 19750
$ws.Range("S25").Value = 988
